$wb = $excel.ActiveWorkbook

# Updates apply identically to the "展览" sheet and the "全部类型" sheet
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 830
    $ws.Range("F6").Value = 12379
    $ws.Range("F12").Value = 914
    $ws.Range("F13").Value = 13620
    $ws.Range("F14").Value = 13829
    $ws.Range("F19").Value = 1035
}
